# Add two new slides ("Hypothèses" and "Sources") at the end of the deck,
# using the existing "Titre et contenu" (Title + Content) layout -- the
# same one already used by slide 2 -- so the new slides pick up matching
# title / content placeholders.

$p = $ppt.ActivePresentation
$titleContentLayout = $p.SlideMaster.CustomLayouts.Item(2)

# --- Slide 3: "Hypothèses" ---------------------------------------------
$s3 = $p.Slides.AddSlide($p.Slides.Count + 1, $titleContentLayout)

$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Hypothèses"

$s3body = $s3.Shapes.Item(2).TextFrame.TextRange
$s3body.Text = "Pas d’éolien off-shore`rPas de solaire au sol -> uniquement toiture"

# --- Slide 4: "Sources" -------------------------------------------------
$s4 = $p.Slides.AddSlide($p.Slides.Count + 1, $titleContentLayout)

$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Sources"

$s4body = $s4.Shapes.Item(2).TextFrame.TextRange
$s4body.Text = "Futurs énergétiques 2050 - RTE`rWorld "
$s4body.InsertAfter("Energy ")
$s4body.InsertAfter("Outlook 2020 - IAE")
